$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the zoom that had crept up to 175% back to the normal 100% level.
$excel.ActiveWindow.Zoom = 100

# Re-insert the "shadda alone" (ّ) row that a prior commit had dropped,
# right before the shadda+vowel combination rows, restoring id 8..15.
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "ّ"

# Renumber the id column for every row pushed down by the insert.
for ($r = 10; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

$ws.Range("A17").Select()
